# Generate Report for Archive
#
# The localization status moved on from "Ready for handoff" to
# "In Translation": update every cell carrying that status (the Overview
# sheet's per-locale status columns, plus each locale sheet's own "Status"
# column), then re-fit the columns that held the old, longer text so their
# width matches the shorter replacement string.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target "character" column width (matches the narrower width the status
# column settles at once it only has to fit "In Translation").
$targetWidth = 12.5

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $targetWidth
$overview.Columns.Item(6).ColumnWidth = $targetWidth

# --- Per-locale sheets: "Status" column (C2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $targetWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $targetWidth
